# Regenerate merged AHB files
# - rename the "_old" / "_new" header suffixes to "_FV2404" / "_FV2410"
# - wrap the data range in an Excel Table (Table1)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename header row (row 1) values.
# ---------------------------------------------------------------------------
$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$headersFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# columns A..J (1..10)
for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2404[$i]
}

# column K (11) stays "diff"
$ws.Cells.Item(1, 11).Value = "diff"

# columns L..U (12..21)
for ($i = 0; $i -lt $headersFV2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2410[$i]
}

# ---------------------------------------------------------------------------
# 2. Turn the data range into an Excel Table ("Table1").
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------------
# 3. Freeze the header row.
# ---------------------------------------------------------------------------
[void]$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
